# Update confusion-matrix values in the "JIND_raw" and "JIND" sheets
# (commit message: "ConfMatrix and pre Nat")

$wb = $excel.ActiveWorkbook

# --- Sheet: JIND_raw ---
$wsRaw = $wb.Worksheets.Item("JIND_raw")
$wsRaw.Range("I3").Value  = 0.015228426395939087
$wsRaw.Range("C8").Value  = 0.001128668171557562
$wsRaw.Range("I9").Value  = 0.9695431472081218
$wsRaw.Range("C10").Value = 0.001128668171557562

# --- Sheet: JIND ---
$wsJind = $wb.Worksheets.Item("JIND")
$wsJind.Range("C3").Value  = 0.9830699774266366
$wsJind.Range("H3").Value  = 0.14285714285714285
$wsJind.Range("I5").Value  = 0.01015228426395939
$wsJind.Range("I9").Value  = 0.9543147208121827
$wsJind.Range("C10").Value = 0.014672686230248307
$wsJind.Range("H10").Value = 0.0

$wb.Save()
